# Refresh the "cryptos" sheet: updated Price (D) and Volume(1h) (E) values
# for every coin, plus re-synced Coin (B) / Link (C) text for rows whose
# ranking shifted up/down since the last snapshot.
#
# Each hashtable below is one changed data row; only the columns that
# actually changed are included (keyed by column letter).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="22.452.52"; E="  +0.21%  " },
    @{ Row=3; D="1.572.96"; E="  +0.57%  " },
    @{ Row=4; E="  +0.01%  " },
    @{ Row=5; E="  +0.06%  " },
    @{ Row=6; D="287.87"; E="  +0.75%  " },
    @{ Row=7; D="0.3709"; E="  +1.85%  " },
    @{ Row=8; D="47.31"; E="  -2.32%  " },
    @{ Row=9; D="0.3318"; E="  -0.47%  " },
    @{ Row=10; D="1.153"; E="  +2.29%  " },
    @{ Row=11; D="0.07504"; E="  +1.27%  " },
    @{ Row=12; E="  +0.03%  " },
    @{ Row=13; D="20.77"; E="  -0.06%  " },
    @{ Row=14; D="5.933"; E="  +0.07%  " },
    @{ Row=15; D="6.919"; E="  +0.47%  " },
    @{ Row=16; D="1.562.46"; E="  -0.07%  " },
    @{ Row=17; D="0.00001116"; E="  +1.08%  " },
    @{ Row=18; D="88.35"; E="  +0.20%  " },
    @{ Row=19; D="0.06728"; E="  +0.53%  " },
    @{ Row=20; D="0.9998"; E="  -0.06%  " },
    @{ Row=21; D="6.383"; E="  +0.44%  " },
    @{ Row=22; D="16.48"; E="  +2.43%  " },
    @{ Row=23; D="11.99"; E="  +0.14%  " },
    @{ Row=24; D="22.446.83"; E="  +0.24%  " },
    @{ Row=25; D="2.387"; E="  -1.17%  " },
    @{ Row=26; D="2.628"; E="  +3.07%  " },
    @{ Row=27; D="150.45"; E="  +0.41%  " },
    @{ Row=28; D="19.60"; E="  +0.96%  " },
    @{ Row=29; D="4.958" },
    @{ Row=30; D="125.16"; E="  +1.63%  " },
    @{ Row=31; D="1.740.90"; E="  +0.21%  " },
    @{ Row=32; D="1.096"; E="  +2.79%  " },
    @{ Row=33; D="6.088"; E="  -0.62%  " },
    @{ Row=34; D="1.985"; E="  -0.38%  " },
    @{ Row=35; D="9.869"; E="  +2.90%  " },
    @{ Row=36; D="0.08328"; E="  +0.97%  " },
    @{ Row=37; D="0.02444"; E="  +2.39%  " },
    @{ Row=38; D="1.308"; E="  +0.13%  " },
    @{ Row=39; D="0.2223"; E="  +0.79%  " },
    @{ Row=40; D="0.06381"; E="  +0.06%  " },
    @{ Row=41; D="5.323"; E="  -0.18%  " },
    @{ Row=42; D="11.38"; E="  +2.51%  " },
    @{ Row=43; D="0.6240"; E="  +2.78%  " },
    @{ Row=44; D="14.04"; E="  +2.61%  " },
    @{ Row=45; B="Decentraland"; C="https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D="0.6064"; E="  +5.57%  " },
    @{ Row=46; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="3.772"; E="  +0.37%  " },
    @{ Row=47; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="2.045"; E="  +1.77%  " },
    @{ Row=48; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="124.83"; E="  +0.12%  " },
    @{ Row=49; B="EOS"; C="https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"; D="1.207"; E="  -0.47%  " },
    @{ Row=50; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.07200" },
    @{ Row=51; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="77.27"; E="  +2.75%  " }
)

foreach ($update in $updates) {
    $row = $update.Row

    foreach ($col in @("B", "C", "D", "E")) {
        if (-not $update.ContainsKey($col)) { continue }

        $value = $update[$col]
        $cell = $ws.Range("$col$row")

        # Every value on this sheet is stored as text - coin names, coin
        # links, price strings such as "287.87", and padded percentages such
        # as "  +0.21%  ". Excel auto-converts plain numeric-looking text
        # (e.g. "287.87", "0.9998") to a real number on assignment, which
        # would corrupt these price cells. Guard against that by writing a
        # leading apostrophe (the classic force-text marker) whenever the
        # raw value is a plain number, then reset the cell style afterwards
        # so the apostrophe trick doesn't leave a stray number format behind.
        if ($value -match '^-?\d+(\.\d+)?$') {
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
